# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the other header cells and filling the data rows with 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new
# header cell (H1) so it keeps the bold/centered/bordered style (s="1"),
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H10 with the value 1 (unstyled, like the other numeric columns).
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
